# Applies the commit's changes to food_table2.xlsx:
#  1. Rename sheets: drop the "Food density " prefix from both tab names.
#  2. Reset the custom row height (16.2pt) back to the sheet's default on a
#     number of data rows in both sheets (these rows had picked up an
#     explicit height from editing and are being auto-fit back to normal).
#  3. Shrink the bottom border row (row 57) on sheet "2017 2018" from 16.8pt
#     down to 15pt.
#  4. Scroll sheet "2017 2018" down so row 22 is at the top of the view, and
#     move the active selection from O51 to O52.

$wb = $excel.ActiveWorkbook

# --- 1. Rename worksheets -------------------------------------------------
$wsYear1 = $wb.Worksheets.Item(1)
$wsYear2 = $wb.Worksheets.Item(2)

$wsYear1.Name = "2015 2016"
$wsYear2.Name = "2017 2018"

# --- 2. Normalize row heights back to default (AutoFit) -------------------
# Sheet "2015 2016": rows whose ht="16.2" reverts to the default row height.
$sheet1Rows = @(29, 35, 47, 49)
foreach ($r in $sheet1Rows) {
    $wsYear1.Rows.Item($r).AutoFit()
}

# Sheet "2017 2018": same reset for the matching set of rows.
$sheet2Rows = @(5, 9, 13, 17, 21, 25, 29, 31, 33, 35, 37, 39, 41, 43, 45, 47, 49, 51, 53, 55)
foreach ($r in $sheet2Rows) {
    $wsYear2.Rows.Item($r).AutoFit()
}

# --- 3. Row 57 on "2017 2018": shrink from 16.8pt to 15pt ------------------
$wsYear2.Rows.Item(57).RowHeight = 15

# --- 4. Update the view/selection on "2017 2018" ---------------------------
$wsYear2.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$wsYear2.Range("O52").Select() | Out-Null
